# Apply the commit: drop the intro/title slides and the closing "Thank you"
# slide, keeping only the five content slides (old positions 4-8), and
# refresh the cached "datetimeFigureOut" field text used by the Date
# placeholders from 09/03/2023 to 18/03/2023 on the slide master and every
# slide layout.

$p = $ppt.ActivePresentation

$oldDate = "09/03/2023"
$newDate = "18/03/2023"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shape = $shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                if ($tf.TextRange.Text -eq $oldDate) {
                    $tf.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master's own Date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every layout's Date placeholder.
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# Delete slides from the end first so earlier indices stay stable:
# old order (1-based): 1=Project name/Hair studio, 2=Silver oak college,
# 3=Subject name: INTERNSHIP OJT, 4=Fresh look haircut services,
# 5=Introduction of the company, 6=History of the haircut,
# 7=Images of the Salon, 8=Haircut services men/women, 9=Thank you
$p.Slides.Item(9).Delete()
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()
$p.Slides.Item(1).Delete()
